# Scen_RES_SHARE_50%_24_7.xlsx — "Add files via upload"
#
# The sheet INS_1 holds a single ~TFM_INS table. Row 6 (the data row) has
# its LimType (column G) flipped from "FX" (fixed) to "LO" (lower bound),
# and the user's last on-screen selection moved from C13 to F15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS_1")

# LimType FX -> LO for the ELC_IND_FIN_DEM / ELC_IND_RES_SUM / FLO_SHAR row.
$ws.Range("G6").Value = "LO"

# Restore the reported active-cell selection.
$ws.Range("F15").Select()
